$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.152.81"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.463.47"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.24"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.15"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +4.47%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.464.06"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.061.03"
$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.38"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.144.21"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.453.51"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.24"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.17"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.559"
$ws.Range("E23").Value = "  +2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.79"
$ws.Range("E24").Value = "  -1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +5.93%  "

$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.51"
$ws.Range("E30").Value = "  +4.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("E32").Value = "  +1.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.60"
$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.35"
$ws.Range("E34").Value = "  +5.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +9.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.85"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("E37").Value = "  +2.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0784"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.60"
$ws.Range("E39").Value = "  +1.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.80"
$ws.Range("E40").Value = "  +8.20%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.90"
$ws.Range("E41").Value = "  +4.64%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.899.69"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0321"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.25"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.782"
$ws.Range("E45").Value = "  +2.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.14"
$ws.Range("E46").Value = "  +10.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "324.00"
$ws.Range("E47").Value = "  +10.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("E49").Value = "  +1.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.881"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.63"
$ws.Range("E51").Value = "  +1.66%  "
